$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A122").Value = "reg-ack-template-part4"
$ws.Range("B122").Value = "Registration Acknowledgement Template - Part 4"
$ws.Range("C122").Value = "eng"
$ws.Range("D122").Value = $true
$ws.Range("E122").Value = "superadmin"
$ws.Range("F122").Value = "now()"

$ws.Range("A123").Value = "reg-ack-template-part4"
$ws.Range("B123").Value = "نموذج شكر التسجيل"
$ws.Range("C123").Value = "ara"
$ws.Range("D123").Value = $true
$ws.Range("E123").Value = "superadmin"
$ws.Range("F123").Value = "now()"

$ws.Range("A124").Value = "reg-ack-template-part4"
$ws.Range("B124").Value = "accusé de réception"
$ws.Range("C124").Value = "fra"
$ws.Range("D124").Value = $true
$ws.Range("E124").Value = "superadmin"
$ws.Range("F124").Value = "now()"

$ws.Range("A125:XFD1048576").Select()
